# Chiffres COVID-19 Valais — add the 2021-03-29 (row 401) data point and
# apply a few corrections to previously-entered "new positive cases" counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to "Nb nouveaux cas positifs" (column C) for earlier days ---
$ws.Range("C336").Value = 88
$ws.Range("C396").Value = 44
$ws.Range("C398").Value = 132
$ws.Range("C399").Value = 88
$ws.Range("C400").Value = 72

# --- New data for row 401 (2021-03-29) ---
$ws.Range("C401").Value = 11   # Nb nouveaux cas positifs
$ws.Range("E401").Value = 7    # Patients COVID-19 aux SI total (y.c. intubes)
$ws.Range("F401").Value = 6    # Patients COVID-19 intubes
$ws.Range("G401").Value = 16   # Patients COVID-19 hospitalises hors SI

# Columns L (Nb nouveaux deces a l'hopital) and M (Nb nouveaux deces
# extra-hospitaliers) are formatted as Text ("@") on this row, inherited
# from the placeholder row. Writing a number straight into a Text-formatted
# cell would store it as text, so briefly switch the cell to a general
# number format, write the numeric 0, then restore the original Text
# format (matches how the rest of the column already stores its values).
$ws.Range("L401").NumberFormat = "general"
$ws.Range("L401").Value = 0
$ws.Range("L401").NumberFormat = "@"

$ws.Range("M401").NumberFormat = "general"
$ws.Range("M401").Value = 0
$ws.Range("M401").NumberFormat = "@"
